$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Travail effectué" text for Jeudi (row 6) in column F
$ws.Range("F6").Value = "Modifications pour relance facture et pdf relance facture"

# Update the active selection to reflect the last edited cell
$ws.Range("F12").Select()
